# Insert a new data row at row 333 (pushes existing rows 333-403 down to
# 334-404, extending the used range from A1:R403 to A1:R404), then populate
# the newly inserted row with the new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("333:333").Insert()

$ws.Cells.Item(333, 1).Value  = 9
$ws.Cells.Item(333, 2).Value  = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(333, 3).Value  = 'Metropolitana'
$ws.Cells.Item(333, 4).Value  = 44785
$ws.Cells.Item(333, 5).Value  = 13
$ws.Cells.Item(333, 6).Value  = 100112052
$ws.Cells.Item(333, 7).Value  = 'Albahaca'
$ws.Cells.Item(333, 8).Value  = 'Sin especificar'
$ws.Cells.Item(333, 9).Value  = 'Primera'
$ws.Cells.Item(333, 10).Value = 80
$ws.Cells.Item(333, 11).Value = 6000
$ws.Cells.Item(333, 12).Value = 6000
$ws.Cells.Item(333, 13).Value = 6000
$ws.Cells.Item(333, 14).Value = '$/paquete'
$ws.Cells.Item(333, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(333, 16).Value = 6000
$ws.Cells.Item(333, 17).Value = 1
$ws.Cells.Item(333, 18).Value = 'Hortaliza'
